$d = $word.ActiveDocument

$replacements = @(
    @{old = "2023-09-24 Sunday"; new = "2023-09-25 Monday"},
    @{old = "52×93=4836"; new = "61×89=5429"},
    @{old = "24×65=1560"; new = "33×70=2310"},
    @{old = "66×47=3102"; new = "78×71=5538"},
    @{old = "93×84=7812"; new = "11×56=616"},
    @{old = "88×97=8536"; new = "90×15=1350"},
    @{old = "61×53=3233"; new = "62×47=2914"},
    @{old = "71×12=852"; new = "25×11=275"},
    @{old = "62×21=1302"; new = "30×59=1770"},
    @{old = "96×87=8352"; new = "37×34=1258"},
    @{old = "51×45=2295"; new = "12×70=840"},
    @{old = "93×68=6324"; new = "15×97=1455"},
    @{old = "27×76=2052"; new = "77×91=7007"},
    @{old = "56×75=4200"; new = "45×82=3690"},
    @{old = "83×37=3071"; new = "27×70=1890"},
    @{old = "42×60=2520"; new = "29×75=2175"},
    @{old = "47×65=3055"; new = "78×84=6552"},
    @{old = "51×46=2346"; new = "14×36=504"},
    @{old = "11×70=770"; new = "73×76=5548"},
    @{old = "67×45=3015"; new = "96×71=6816"},
    @{old = "56×33=1848"; new = "77×91=7007"},
    @{old = "64×56=3584"; new = "21×90=1890"},
    @{old = "57×77=4389"; new = "45×50=2250"},
    @{old = "27×75=2025"; new = "35×42=1470"},
    @{old = "68×59=4012"; new = "72×57=4104"},
    @{old = "65×28=1820"; new = "71×87=6177"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
